$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 9 (August) label and values
$ws.Range("A9").Value = "August (through 08-15)"
$ws.Range("B9").Value = 15
$ws.Range("C9").Value = 35
$ws.Range("D9").Value = 34
$ws.Range("E9").Value = 24
$ws.Range("F9").Value = 20
$ws.Range("G9").Value = 93
$ws.Range("H9").Value = 85

# Update row 10 (Total) values
$ws.Range("B10").Value = 177
$ws.Range("C10").Value = 337
$ws.Range("D10").Value = 499
$ws.Range("E10").Value = 449
$ws.Range("F10").Value = 324
$ws.Range("G10").Value = 714
$ws.Range("H10").Value = 1000
